# Apply the "ambios reporte excel correcciones" changes to the first sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) is unchanged ---

# --- Row 2: pelusa / perro / blanco / macho / Publicado ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "pelusa"
$ws.Range("C2").Value = 5
$ws.Range("D2").Value = "perro"
$ws.Range("E2").Value = "blanco"
$ws.Range("F2").Value = "macho"
$ws.Range("G2").Value = "Publicado"
$ws.Range("H2").Value = "/assets/recibidas/1.jpg"
$ws.Range("I2").Value = '{"lat":-34.60975,"lng":-58.428904}'

# --- Row 3: Truman / perro / rubio / macho / Publicado ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Truman"
$ws.Range("C3").Value = 4
$ws.Range("D3").Value = "perro"
$ws.Range("E3").Value = "rubio"
$ws.Range("F3").Value = "macho"
$ws.Range("G3").Value = "Publicado"
$ws.Range("H3").Value = ""
$ws.Range("I3").Value = '{"lat":-34.613464,"lng":-58.428317}'

# --- Row 4: Truman / perro / rubio / macho / Publicado ---
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Truman"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = "perro"
$ws.Range("E4").Value = "rubio"
$ws.Range("F4").Value = "macho"
$ws.Range("G4").Value = "Publicado"
$ws.Range("H4").Value = ""
$ws.Range("I4").Value = '{"lat":-34.613464,"lng":-58.428317}'

# --- Row 5: Truman / perro / rubio / macho / Publicado ---
$ws.Range("A5").Value = 2
$ws.Range("B5").Value = "Truman"
$ws.Range("C5").Value = 4
$ws.Range("D5").Value = "perro"
$ws.Range("E5").Value = "rubio"
$ws.Range("F5").Value = "macho"
$ws.Range("G5").Value = "Publicado"
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = '{"lat":-34.613464,"lng":-58.428317}'

# --- Row 6: Truman / perro / rubio / macho / Publicado ---
$ws.Range("A6").Value = 2
$ws.Range("B6").Value = "Truman"
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = "perro"
$ws.Range("E6").Value = "rubio"
$ws.Range("F6").Value = "macho"
$ws.Range("G6").Value = "Publicado"
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = '{"lat":-34.613464,"lng":-58.428317}'

# --- Row 7: Truman / perro / rubio / macho / Publicado ---
$ws.Range("A7").Value = 2
$ws.Range("B7").Value = "Truman"
$ws.Range("C7").Value = 4
$ws.Range("D7").Value = "perro"
$ws.Range("E7").Value = "rubio"
$ws.Range("F7").Value = "macho"
$ws.Range("G7").Value = "Publicado"
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = '{"lat":-34.613464,"lng":-58.428317}'

# --- Row 8: Truman / perro / rubio / macho / Publicado, photo 3 ---
$ws.Range("A8").Value = 3
$ws.Range("B8").Value = "Truman"
$ws.Range("C8").Value = 4
$ws.Range("D8").Value = "perro"
$ws.Range("E8").Value = "rubio"
$ws.Range("F8").Value = "macho"
$ws.Range("G8").Value = "Publicado"
$ws.Range("H8").Value = "/assets/recibidas/3.jpg"
$ws.Range("I8").Value = ""

# --- Row 9: Truman / perro / rubio / macho / Publicado, photo 4 ---
$ws.Range("A9").Value = 4
$ws.Range("B9").Value = "Truman"
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "perro"
$ws.Range("E9").Value = "rubio"
$ws.Range("F9").Value = "macho"
$ws.Range("G9").Value = "Publicado"
$ws.Range("H9").Value = "/assets/recibidas/4.jpg"
$ws.Range("I9").Value = ""
